$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Pawlotsky, 2016
$ws.Range("A5").Value = "Pawlotsky_2016"
$ws.Range("B5").Value = "Pawlotsky, 2016"
$ws.Range("C5").Value = "Pawlotsky, J-M"
$ws.Range("D5").Value = 2016
$ws.Range("E5").Value = "Hepatitis C Virus Resistance to Direct-Acting Antiviral Drugs in Interferon-Free Regimens"
$ws.Range("F5").Value = "Gastroenterology"
$ws.Range("G3").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G5").Value = "151(1):70-86"
$ws.Range("H5").Value = "http://www.sciencedirect.com/science/article/pii/S0016508516300555"

# Row 6: Sarrazin, 2016
$ws.Range("A6").Value = "Sarrazin_2016"
$ws.Range("B6").Value = "Sarrazin, 2016"
$ws.Range("C6").Value = "Sarrazin, C"
$ws.Range("D6").Value = 2016
$ws.Range("E6").Value = "The importance of resistance to direct antiviral drugs in HCV infection in clinical practice"
$ws.Range("F6").Value = "Journal of Hepatology"
$ws.Range("G3").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G6").Value = "64(2):486-504"
$ws.Range("H6").Value = "http://www.sciencedirect.com/science/article/pii/S0168827815006297"

# Row 7: Lontok et al., 2015
$ws.Range("A7").Value = "Lontok_et_al_2015"
$ws.Range("B7").Value = "Lontok et al., 2015"
$ws.Range("C4").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = "Lontok E, Harrington P, Howe A, Kieffer T, Lennerstrand J, Lenz O, McPhee F, Mo H, Parkin N, Pilot-Matias T, Miller V"
$ws.Range("D7").Value = 2015
$ws.Range("E7").Value = "Hepatitis C virus drug resistance-associated substitutions: State of the art summary"
$ws.Range("F7").Value = "Hepatology"
$ws.Range("G3").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").Value = "62(5):1623-32"
$ws.Range("H7").Value = "http://onlinelibrary.wiley.com/doi/10.1002/hep.27934/abstract;jsessionid=657E6E3C85196C806137825E68AA9660.f02t01"

$ws.Range("G10").Select()

Write-Host "done"
